# issue #5: add legislator_id, name, date into dataframe
#
# The per-legislator scraper output gets three new trailing columns on the
# stock ("股票") sheet: the filing date, the legislator's name and their
# numeric id. Every existing data row gets the same three values appended;
# the header row gets three new column headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

$legislatorDate = "2013-12-26"
$legislatorName = "丁守中"
$legislatorId = 515

# How many data rows already exist below the header (row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 11 }

# --- Header row ------------------------------------------------------
# Reuse the header formatting (bold, bordered, centered) already used by
# the existing header cells by copying it from the last header cell (G1).
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows ---------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    # Match the plain data-row formatting used by the rest of the row.
    $ws.Range("G$r").Copy()
    $ws.Range("H$r`:J$r").PasteSpecial(-4122)

    # The date column would otherwise be auto-recognised as a real date by
    # Excel and stored as a date serial number; force it to stay a literal
    # text value (as the scraper wrote it) by pre-formatting as Text, then
    # restoring the plain formatting used by the rest of the row.
    $ws.Range("H$r").NumberFormat = "@"
    $ws.Range("H$r").Value = $legislatorDate
    $ws.Range("G$r").Copy()
    $ws.Range("H$r").PasteSpecial(-4122)

    $ws.Range("I$r").Value = $legislatorName
    $ws.Range("J$r").Value = $legislatorId
}

$excel.CutCopyMode = $false
